$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.619.17'
$ws.Range("E2").Value = '  +3.08%  '
$ws.Range("D3").Value = '''2.433.91'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''563.99'
$ws.Range("E5").Value = '  +2.22%  '
$ws.Range("D6").Value = '''166.43'
$ws.Range("E6").Value = '  +4.30%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '''0.513'
$ws.Range("E8").Value = '  +1.18%  '
$ws.Range("D9").Value = '''0.169'
$ws.Range("E9").Value = '  +7.01%  '
$ws.Range("D10").Value = '''2.432.66'
$ws.Range("E10").Value = '  +0.40%  '
$ws.Range("E11").Value = '  -2.11%  '
$ws.Range("D12").Value = '''0.335'
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").Value = '''4.67'
$ws.Range("E13").Value = '  -2.11%  '
$ws.Range("D14").Value = '''0.0000178'
$ws.Range("E14").Value = '  +4.85%  '
$ws.Range("D15").Value = '''69.463.51'
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("D16").Value = '''2.880.73'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").Value = '''23.97'
$ws.Range("E17").Value = '  +4.41%  '
$ws.Range("D18").Value = '''2.430.34'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = '''10.81'
$ws.Range("E19").Value = '  +4.49%  '
$ws.Range("D20").Value = '''341.92'
$ws.Range("E20").Value = '  +4.04%  '
$ws.Range("D21").Value = '''7.14'
$ws.Range("E21").Value = '  +4.30%  '
$ws.Range("D22").Value = '''3.89'
$ws.Range("E22").Value = '  +2.79%  '
$ws.Range("E23").Value = '  +6.74%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '''65.92'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").Value = '''3.82'
$ws.Range("E26").Value = '  +5.88%  '
$ws.Range("D27").Value = '''8.52'
$ws.Range("E27").Value = '  +5.71%  '
$ws.Range("D28").Value = '''2.554.83'
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '''0.0₃0851'
$ws.Range("E30").Value = '  +6.05%  '
$ws.Range("D31").Value = '''7.39'
$ws.Range("E31").Value = '  +5.40%  '
$ws.Range("D32").Value = '''1.24'
$ws.Range("E32").Value = '  +10.27%  '
$ws.Range("D33").Value = '''452.39'
$ws.Range("E33").Value = '  +8.75%  '
$ws.Range("D34").Value = '''0.999'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("D36").Value = '''157.42'
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("D37").Value = '''19.13'
$ws.Range("E37").Value = '  +1.08%  '
$ws.Range("D38").Value = '''0.111'
$ws.Range("E38").Value = '  +5.53%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '''18.23'
$ws.Range("E40").Value = '  +2.58%  '
$ws.Range("D41").Value = '''0.302'
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '''4.40'
$ws.Range("E42").Value = '  +3.67%  '
$ws.Range("D43").Value = '''1.52'
$ws.Range("E43").Value = '  +4.55%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''37.85'
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("D45").Value = '''1.09'
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("E46").Value = '  +5.22%  '
$ws.Range("D47").Value = '''135.02'
$ws.Range("E47").Value = '  +3.91%  '
$ws.Range("D48").Value = '''3.40'
$ws.Range("E48").Value = '  +2.56%  '
$ws.Range("D49").Value = '''0.0725'
$ws.Range("E49").Value = '  +2.63%  '
$ws.Range("D50").Value = '''0.489'
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("D51").Value = '''0.563'
$ws.Range("E51").Value = '  +1.82%  '
